# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output data (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6650
$ws1.Range("F5").Value = 66
$ws1.Range("F12").Value = 169
$ws1.Range("F16").Value = 3316
$ws1.Range("F19").Value = 1957
$ws1.Range("F20").Value = 73

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6650
$ws4.Range("F5").Value = 66
$ws4.Range("F13").Value = 169
$ws4.Range("F17").Value = 3316
$ws4.Range("F20").Value = 1957
$ws4.Range("F21").Value = 73
